# edit.ps1
# Adds a new "2022-Q1" worksheet (fund-holdings detail, inserted right
# before "总计") and prepends a matching summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet, positioned immediately before "总计"
# ---------------------------------------------------------------------
$totalBefore = $wb.Worksheets.Item("总计")
$new = $wb.Worksheets.Add($totalBefore)
$new.Name = "2022-Q1"

# NOTE: sheet handles in this engine are positional, not stable object
# references -- after Add() shifts indices, re-resolve "总计" by name
# rather than reusing the handle obtained before the insert.
$total = $wb.Worksheets.Item("总计")

# Reuse formatting from an existing quarter sheet's header / index column
$src = $wb.Worksheets.Item("2021-Q4")

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $new.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$src.Range("B1:H1").Copy()
$new.Range("B1:H1").PasteSpecial(-4122)

$data = @"
006567	中泰星元价值优选灵活配置混合	44.13	81.95	4.75	2.0962	6
013776	中泰兴为价值精选混合A	20.31	85.34	4.21	0.8551	8
006624	中泰玉衡价值优选混合	17.75	81.95	4.69	0.8325	6
003293	易方达科瑞灵活配置混合	34.67	78.17	1.95	0.6761	10
013777	中泰兴为价值精选混合C	8.71	85.34	4.21	0.3667	8
400003	东方精选混合	10.47	82.20	2.87	0.3005	10
519021	国泰金鼎价值混合	6.53	90.06	4.11	0.2684	8
519150	新华优选消费混合	4.37	93.57	6.03	0.2635	2
161609	融通动力先锋混合	7.12	80.93	2.96	0.2108	9
009181	浙商智多兴稳健回报一年持有期混合A	17.79	32.21	1.02	0.1815	9
009205	兴银丰运稳益回报混合A	7.11	32.99	1.37	0.0974	9
519093	新华钻石品质企业混合	1.36	92.58	6.45	0.0877	2
001152	融通新区域新经济灵活配置混合	2.81	80.98	2.97	0.0835	8
009182	浙商智多兴稳健回报一年持有期混合C	6.80	32.21	1.02	0.0694	9
000963	兴业多策略灵活配置混合	2.07	75.34	3.31	0.0685	7
009206	兴银丰运稳益回报混合C	4.32	32.99	1.37	0.0592	9
001004	新华稳健回报灵活配置混合	0.77	93.15	6.47	0.0498	2
020033	国泰民安增利债券A	2.12	48.78	2.29	0.0485	6
012461	西藏东财国证龙头家电指数型发起式证券投资基金A	1.01	94.99	3.36	0.0339	7
020034	国泰民安增利债券C	1.20	48.78	2.29	0.0275	6
001626	国泰央企改革股票	0.58	90.32	4.17	0.0242	8
011775	格林鑫悦一年持有期混合型证券投资基金A	2.64	33.74	0.78	0.0206	10
009128	明亚价值长青混合A	0.38	49.48	4.30	0.0163	2
012462	西藏东财国证龙头家电指数型发起式证券投资基金C	0.42	94.99	3.36	0.0141	7
010777	浙商智选家居股票A	0.15	90.92	6.54	0.0098	7
003980	中银证券瑞益灵活配置混合A	0.93	33.74	0.78	0.0073	10
010778	浙商智选家居股票C	0.03	90.92	6.54	0.0020	7
012268	浙商智多享稳健混合型发起式证券投资基金A	0.17	22.32	0.90	0.0015	4
012269	浙商智多享稳健混合型发起式证券投资基金C	0.01	22.32	0.90	0.0001	4
009129	明亚价值长青混合C	0.00	49.48	4.30	0	2
"@

$rows = $data -split "`n"
$r = 2
foreach ($line in $rows) {
    $cols = $line -split "`t"
    $code = $cols[0]
    $name = $cols[1]
    $scale = $cols[2]
    $position = $cols[3]
    $ratio = $cols[4]
    $mktval = $cols[5]
    $rank = [int]$cols[6]

    $new.Cells.Item($r, 1).Value = $r - 2
    $new.Cells.Item($r, 2).Value = "'" + $code
    $new.Cells.Item($r, 3).Value = "'" + $name
    $new.Cells.Item($r, 4).Value = "'" + $scale
    $new.Cells.Item($r, 5).Value = "'" + $position
    $new.Cells.Item($r, 6).Value = "'" + $ratio
    if ($mktval -eq "0") {
        $new.Cells.Item($r, 7).Value = 0
    } else {
        $new.Cells.Item($r, 7).Value = "'" + $mktval
    }
    $new.Cells.Item($r, 8).Value = $rank

    # Drop the auto "quote-prefix" style the text-forcing apostrophe
    # above leaves behind, so B:G stay on the default ("Normal") style,
    # matching the quarter-sheets' plain un-styled data cells.
    $new.Range("B" + $r + ":G" + $r).Style = "Normal"

    $r = $r + 1
}

$src.Range("A2").Copy()
$new.Range("A2:A31").PasteSpecial(-4122)
for ($i = 2; $i -le 31; $i++) {
    $new.Cells.Item($i, 1).Value = $i - 2
}

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" summary row to "总计", shifting old rows down
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 30
$total.Range("D2").Value = 6.77

for ($r = 3; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

Write-Host "2022-Q1 sheet added and 总计 updated"
